{"js": "// Map of old multiplication expressions -> new ones (document order).\nconst replacements = [\n  [\"844\u00d75=\", \"880\u00d76=\"],\n  [\"375\u00d76=\", \"205\u00d72=\"],\n  [\"234\u00d78=\", \"455\u00d75=\"],\n  [\"876\u00d76=\", \"162\u00d79=\"],\n  [\"984\u00d72=\", \"425\u00d75=\"],\n  [\"957\u00d75=\", \"652\u00d77=\"],\n  [\"114\u00d72=\", \"203\u00d78=\"],\n  [\"991\u00d77=\", \"974\u00d79=\"],\n  [\"710\u00d72=\", \"130\u00d72=\"],\n  [\"539\u00d75=\", \"251\u00d79=\"],\n  [\"589\u00d73=\", \"856\u00d76=\"],\n  [\"380\u00d73=\", \"597\u00d78=\"],\n  [\"162\u00d77=\", \"588\u00d78=\"],\n  [\"520\u00d72=\", \"650\u00d72=\"],\n  [\"315\u00d76=\", \"571\u00d74=\"],\n  [\"638\u00d72=\", \"840\u00d73=\"],\n  [\"367\u00d79=\", \"117\u00d74=\"],\n  [\"642\u00d77=\", \"337\u00d78=\"],\n  [\"493\u00d77=\", \"696\u00d77=\"],\n  [\"490\u00d79=\", \"130\u00d72=\"],\n  [\"534\u00d78=\", \"861\u00d77=\"],\n  [\"747\u00d79=\", \"584\u00d73=\"],\n  [\"210\u00d77=\", \"190\u00d73=\"],\n  [\"499\u00d79=\", \"148\u00d74=\"],\n  [\"234\u00d79=\", \"180\u00d73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Search text not found: ${oldText}`);\n  }\n\n  // Replace only the first match so repeated \"before\" values (if any) are\n  // handled in document order across successive iterations.\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Update the multiplication problems in the practice table.\n# Each \"before\" expression is unique in the document, so a plain\n# Find/Replace (one hit each) reproduces the diff exactly.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"844\u00d75=\"; New = \"880\u00d76=\" },\n    @{ Old = \"375\u00d76=\"; New = \"205\u00d72=\" },\n    @{ Old = \"234\u00d78=\"; New = \"455\u00d75=\" },\n    @{ Old = \"876\u00d76=\"; New = \"162\u00d79=\" },\n    @{ Old = \"984\u00d72=\"; New = \"425\u00d75=\" },\n    @{ Old = \"957\u00d75=\"; New = \"652\u00d77=\" },\n    @{ Old = \"114\u00d72=\"; New = \"203\u00d78=\" },\n    @{ Old = \"991\u00d77=\"; New = \"974\u00d79=\" },\n    @{ Old = \"710\u00d72=\"; New = \"130\u00d72=\" },\n    @{ Old = \"539\u00d75=\"; New = \"251\u00d79=\" },\n    @{ Old = \"589\u00d73=\"; New = \"856\u00d76=\" },\n    @{ Old = \"380\u00d73=\"; New = \"597\u00d78=\" },\n    @{ Old = \"162\u00d77=\"; New = \"588\u00d78=\" },\n    @{ Old = \"520\u00d72=\"; New = \"650\u00d72=\" },\n    @{ Old = \"315\u00d76=\"; New = \"571\u00d74=\" },\n    @{ Old = \"638\u00d72=\"; New = \"840\u00d73=\" },\n    @{ Old = \"367\u00d79=\"; New = \"117\u00d74=\" },\n    @{ Old = \"642\u00d77=\"; New = \"337\u00d78=\" },\n    @{ Old = \"493\u00d77=\"; New = \"696\u00d77=\" },\n    @{ Old = \"490\u00d79=\"; New = \"130\u00d72=\" },\n    @{ Old = \"534\u00d78=\"; New = \"861\u00d77=\" },\n    @{ Old = \"747\u00d79=\"; New = \"584\u00d73=\" },\n    @{ Old = \"210\u00d77=\"; New = \"190\u00d73=\" },\n    @{ Old = \"499\u00d79=\"; New = \"148\u00d74=\" },\n    @{ Old = \"234\u00d79=\"; New = \"180\u00d73=\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Forward = $true\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Wrap = 1\n    $find.Execute([ref]$pair.Old, [ref]$true, [ref]$false, [ref]$false, [ref]$false, [ref]$false, [ref]$true, [ref]1, [ref]$false, [ref]$pair.New, [ref]2)\n}\n"}
